$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 221.78572
$ws.Range("I33").Value = 196.45454
$ws.Range("K33").Value = 196.45454
$ws.Range("M33").Value = 32.54545999999999

# Row 53
$ws.Range("H53").Value = 517.2308
$ws.Range("I53").Value = 154.57143
$ws.Range("J53").Value = 940.3333
$ws.Range("K53").Value = 154.57143
$ws.Range("L53").Value = 940.3333
$ws.Range("M53").Value = 482.42857
$ws.Range("N53").Value = -2214.3333

# Row 62
$ws.Range("H62").Value = 3873.5

# Row 65
$ws.Range("H65").Value = 3873.5

# Row 98
$ws.Range("H98").Value = 988.6
$ws.Range("I98").Value = 989
$ws.Range("J98").Value = 988
$ws.Range("K98").Value = 989
$ws.Range("L98").Value = 988
$ws.Range("M98").Value = 509
$ws.Range("N98").Value = -3984

# Row 111
$ws.Range("H111").Value = 1993.3334
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()

# Row 122
$ws.Range("H122").Value = 988.6
$ws.Range("I122").Value = 989
$ws.Range("J122").Value = 988
$ws.Range("K122").Value = 2967
$ws.Range("L122").Value = 2964
$ws.Range("M122").Value = -517
$ws.Range("N122").Value = -7864

# Row 134
$ws.Range("H134").Value = 80000
$ws.Range("J134").Value = 80000
$ws.Range("L134").Value = 80000
$ws.Range("N134").Value = -90140

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3218.7827
$ws.Range("I45").Value = 1610.8572
$ws.Range("J45").Value = 3922.25
$ws.Range("K45").Value = 1610.8572
$ws.Range("L45").Value = 3922.25
$ws.Range("M45").Value = -1233.8572
$ws.Range("N45").Value = -4676.25

# Row 61
$ws.Range("H61").Value = 675.36365
$ws.Range("I61").Value = 442.9
$ws.Range("K61").Value = 442.9
$ws.Range("M61").Value = -230.9

# Row 97
$ws.Range("H97").Value = 1428
$ws.Range("I97").Value = 908
$ws.Range("K97").Value = 908
$ws.Range("M97").Value = -412

# Row 136
$ws.Range("H136").Value = 675.36365
$ws.Range("I136").Value = 442.9
$ws.Range("K136").Value = 1328.7
$ws.Range("M136").Value = 1221.3

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 4333.3335
$ws.Range("I20").Value = 4000
$ws.Range("J20").Value = 4500
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 4500
$ws.Range("M20").Value = -3753
$ws.Range("N20").Value = -4994

# Row 86
$ws.Range("H86").Value = 3376.6924
$ws.Range("I86").Value = 3114.8572
$ws.Range("K86").Value = 3114.8572
$ws.Range("M86").Value = -1991.8572

# Row 89
$ws.Range("H89").Value = 3376.6924
$ws.Range("I89").Value = 3114.8572
$ws.Range("K89").Value = 15574.286
$ws.Range("M89").Value = -9958.286

# Row 95
$ws.Range("H95").Value = 28749.75
$ws.Range("J95").Value = 28749.75
$ws.Range("L95").Value = 28749.75
$ws.Range("N95").Value = -34241.75

# Row 97
$ws.Range("H97").Value = 10745.6
$ws.Range("I97").Value = 9682
$ws.Range("J97").Value = 15000
$ws.Range("K97").Value = 9682
$ws.Range("L97").Value = 15000
$ws.Range("M97").Value = -8691
$ws.Range("N97").Value = -16982

# Row 134
$ws.Range("H134").Value = 1197.45
$ws.Range("I134").Value = 1197.45
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3592.35
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1057.35
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 505.5
$ws.Range("I16").Value = 505.5
$ws.Range("K16").Value = 505.5
$ws.Range("M16").Value = -218.5

# Row 58
$ws.Range("H58").Value = 1028
$ws.Range("I58").Value = 1037.3334
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 1037.3334
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = -834.3334
$ws.Range("N58").Value = -1406

# Row 86
$ws.Range("H86").Value = 7950686.5
$ws.Range("I86").Value = 9937358
$ws.Range("K86").Value = 9937358
$ws.Range("M86").Value = -9936235

# Row 89
$ws.Range("H89").Value = 7950686.5
$ws.Range("I89").Value = 9937358
$ws.Range("K89").Value = 49686790
$ws.Range("M89").Value = -49681174

# Row 93
$ws.Range("H93").Value = 2966.6667
$ws.Range("I93").Value = 2966.6667
$ws.Range("K93").Value = 2966.6667
$ws.Range("M93").Value = -1094.6667

# Row 94
$ws.Range("H94").Value = 4988.625
$ws.Range("I94").Value = 5250
$ws.Range("K94").Value = 5250
$ws.Range("M94").Value = -4799

# Row 99
$ws.Range("H99").Value = 6323.3335
$ws.Range("I99").Value = 6488.3
$ws.Range("K99").Value = 6488.3
$ws.Range("M99").Value = -4990.3

# Row 105
$ws.Range("H105").Value = 13641.25
$ws.Range("I105").Value = 17521.666
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 17521.666
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -15774.666
$ws.Range("N105").Value = -5494

# Row 113
$ws.Range("H113").Value = 505.5
$ws.Range("I113").Value = 505.5
$ws.Range("K113").Value = 505.5
$ws.Range("M113").Value = 1664.5

# Row 122
$ws.Range("H122").Value = 2720.7144
$ws.Range("I122").Value = 2739.1667
$ws.Range("K122").Value = 8217.500100000001
$ws.Range("M122").Value = -5767.500100000001

# Row 126
$ws.Range("H126").Value = 6323.3335
$ws.Range("I126").Value = 6488.3
$ws.Range("K126").Value = 19464.9
$ws.Range("M126").Value = -16994.9

# Row 136
$ws.Range("H136").Value = 1028
$ws.Range("I136").Value = 1037.3334
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 3112.0002
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -562.0001999999999
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 460.94116
$ws.Range("I4").Value = 489.625
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 1468.875
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = -1356.875
$ws.Range("N4").Value = -230

# Row 26
$ws.Range("H26").Value = 1066.6666
$ws.Range("I26").Value = 1066.6666
$ws.Range("K26").Value = 3199.9998
$ws.Range("M26").Value = -2911.9998

# Row 93
$ws.Range("H93").Value = 250
$ws.Range("I93").Value = 250
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 750
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 1122
$ws.Range("N93").ClearContents()

# Row 114
$ws.Range("H114").Value = 397
$ws.Range("I114").Value = 227
$ws.Range("J114").Value = 453.66666
$ws.Range("K114").Value = 681
$ws.Range("L114").Value = 1360.99998
$ws.Range("M114").Value = 2573
$ws.Range("N114").Value = -7868.999980000001

# Row 117
$ws.Range("H117").Value = 3399.25
$ws.Range("J117").Value = 3399.25
$ws.Range("L117").Value = 10197.75
$ws.Range("N117").Value = -17081.75

# Row 121
$ws.Range("H121").Value = 680.75
$ws.Range("I121").Value = 487.5
$ws.Range("J121").Value = 874
$ws.Range("K121").Value = 1462.5
$ws.Range("L121").Value = 2622
$ws.Range("M121").Value = -152.5
$ws.Range("N121").Value = -5242

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 25000
$ws.Range("J15").Value = 25000
$ws.Range("L15").Value = 25000
$ws.Range("N15").Value = -25576

# Row 81
$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996

# Row 84
$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984

# Row 102
$ws.Range("H102").Value = 2972
$ws.Range("I102").Value = 2720.6667
$ws.Range("J102").Value = 3349
$ws.Range("K102").Value = 2720.6667
$ws.Range("L102").Value = 3349
$ws.Range("M102").Value = -1098.6667
$ws.Range("N102").Value = -6593

# Row 122
$ws.Range("H122").Value = 1673.8
$ws.Range("I122").Value = 1425
$ws.Range("J122").Value = 2669
$ws.Range("K122").Value = 4275
$ws.Range("L122").Value = 8007
$ws.Range("M122").Value = -1825
$ws.Range("N122").Value = -12907

# Row 132
$ws.Range("H132").Value = 4999.75
$ws.Range("I132").Value = 4999.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 14999.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12469.25
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2579.6667
$ws.Range("I40").Value = 2579.6667
$ws.Range("K40").Value = 2579.6667
$ws.Range("M40").Value = -2443.6667

# Row 46
$ws.Range("H46").Value = 4060.3914
$ws.Range("I46").Value = 3999.1667
$ws.Range("K46").Value = 3999.1667
$ws.Range("M46").Value = -3811.1667

# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

# Row 55
$ws.Range("H55").Value = 553.625
$ws.Range("I55").Value = 131
$ws.Range("K55").Value = 131
$ws.Range("M55").Value = 42

# Row 122
$ws.Range("H122").Value = 4233.3335
$ws.Range("I122").Value = 4250
$ws.Range("J122").Value = 4200
$ws.Range("K122").Value = 12750
$ws.Range("L122").Value = 12600
$ws.Range("M122").Value = -10300
$ws.Range("N122").Value = -17500

# Row 132
$ws.Range("H132").Value = 3255.3333
$ws.Range("I132").Value = 3255.3333
$ws.Range("K132").Value = 9765.999899999999
$ws.Range("M132").Value = -7235.999899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2554.111
$ws.Range("I122").Value = 2284
$ws.Range("J122").Value = 3499.5
$ws.Range("K122").Value = 6852
$ws.Range("L122").Value = 10498.5
$ws.Range("M122").Value = -4402
$ws.Range("N122").Value = -15398.5

# Row 126
$ws.Range("H126").Value = 4193.8
$ws.Range("I126").Value = 4489.6665
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 13468.9995
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -10998.9995
$ws.Range("N126").Value = -16190

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
